$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level metadata -------------------------------------------------
# (x15ac:absPath is internal Excel last-saved-path metadata; it is not exposed
# through the Excel object model, so it cannot be set from COM automation and is
# left untouched here.)

# --- Row 4: TreePoseVideoPath / hyperlink to the video -----------------------
$ws.Range("B4").Value = "\\Videos\\TreePose.mp4"
$ws.Hyperlinks.Add($ws.Range("B4"), "Videos\TreePose.mp4")
$ws.Range("A4").Value = "TreePoseVideoPath"

# --- Row 5: TreePoseVideoDescription / long pose description ----------------
$ws.Range("A5").Value = "TreePoseVideoDescription"

$description = @'
Tree Pose or Vrksasana
 • From standing pose, shift your weight on to the left leg. Ground your left foot down into the floor, pull up your knee cap and thigh.
• Gaze at a point on eye level, to help you balance.
• Bend your right knee, reach down with your right hand and clasp your ankle.
• With help of the hand, place your right foot on the inside of the left leg, if its available to you on the inner thigh, with your heel up high. Lengthen your tailbone toward the floor.
• If that is not available to you, place the inside of the foot on the inside of the ankle or calf. Avoid the inside of the knee.
• Turn the knee out to the side, press your foot against the inner thigh, and the inner thigh back into the foot.
• Bring your hands in front of your heart in prayer position, or you can lift your arms up to the ceiling. hands shoulderwidth apart.
• With arms lifted, lift from the back body.
• Stay in this pose for about 30 sec or 1 min.
• On an outbreath lower the leg and arms down and stand straight, repeat on the other side.
BENEFITS:
• Improves balance
• Opens the hips
• Strengthens the ankles, legs and spine
• Lengthens the spine
• Improves focus/ concentration
Things to keep in Mind:
• Make sure you keep both arm extended.
• Try to avoid collapsing the side body - keep your chest expanded.

'@

$ws.Range("B5").Value = $description
$ws.Range("B5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 409.5

# --- Selection (best-effort; this engine does not persist scroll position ----
# unless freeze panes/split are active, so "topLeftCell" cannot be reproduced) -
$ws.Range("A5").Select()

Write-Host "done"
